$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Program manager with over 10 years ..." -> "... over 10+ years ..."
#
# The target OOXML keeps this as three separate runs (all sharing the
# same rPr) with the relocated "_GoBack" bookmark sitting between the
# "+" run and the " years ..." run:
#
#   run1: "Program manager with over 10"
#   run2: "+"
#   <bookmarkStart/End w:name="_GoBack">
#   run3: " years of extensive experience ... technology programs."
#
# Plain Range.InsertAfter()/Range.Text= on this paragraph causes the
# host engine to coalesce *every* same-formatted run in the paragraph
# (the many one-word runs later in the sentence), which we must not
# disturb. Inserting through a Bookmark's .Range avoids that
# whole-paragraph coalesce; a scoped Bookmarks.Add/Delete pair is then
# used purely to force the run boundaries we need, which also doesn't
# trigger the coalesce.
# ------------------------------------------------------------------

# 1) Jump to the point right after "...over 10" and drop in a
#    throw-away bookmark there so we have a handle on that exact spot.
$anchor = $d.Content
$ok = $anchor.Find.Execute("Program manager with over 10", $true, $false,
                            $false, $false, $false, $true, 1, $false,
                            "", 0)
if (-not $ok) { throw "Could not find 'Program manager with over 10'" }
$anchor.Collapse(0)
$d.Bookmarks.Add("ZZ_PLUS_POS", $anchor)

# 2) Insert the literal "+" through the bookmark's Range (not the
#    Range returned by Find) - this keeps the edit local to this one
#    paragraph instead of re-flattening every run in it.
$plusSpot = $d.Bookmarks("ZZ_PLUS_POS").Range
$plusSpot.InsertBefore("+")
$d.Bookmarks("ZZ_PLUS_POS").Delete()

# 3) "+" is currently glued onto the end of the "...over 10" run
#    (same formatting). Re-find "...over 10" and plant/drop another
#    scoped bookmark right there to force the run split between
#    "...over 10" and "+".
$splitBefore = $d.Content
$ok = $splitBefore.Find.Execute("Program manager with over 10", $true, $false,
                                 $false, $false, $false, $true, 1, $false,
                                 "", 0)
if (-not $ok) { throw "Could not re-find 'Program manager with over 10'" }
$splitBefore.Collapse(0)
$d.Bookmarks.Add("ZZ_SPLIT_BEFORE_PLUS", $splitBefore)
$d.Bookmarks("ZZ_SPLIT_BEFORE_PLUS").Delete()

# 4) Split "+" from the " years ..." text that follows it, and plant
#    the relocated "_GoBack" bookmark right at that boundary. Word
#    keeps only one "_GoBack" bookmark, so this also removes the one
#    that used to sit right after "Bit Bucket".
$splitAfter = $d.Content
$ok = $splitAfter.Find.Execute("Program manager with over 10+", $true, $false,
                                $false, $false, $false, $true, 1, $false,
                                "", 0)
if (-not $ok) { throw "Could not find 'Program manager with over 10+'" }
$splitAfter.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitAfter)
